# Refactor currency conversion sheet: split "foreign_amount" into explicit
# source_amount / target_amount (+ target_fees) columns.
#
# Before: date | foreign_amount | source_fees | source_currency | target_currency | comment
# After:  date | source_amount  | source_fees | source_currency | target_amount | target_fees | target_currency | comment

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_conversions")

# Make room for the two new columns (target_amount, target_fees) by shifting
# the existing target_currency/comment columns two places to the right.
$ws.Range("E1:F1").Insert(-4161)

# The insert grows the sheet's total column count by 2; shrink it back down
# by removing two blank columns far out of the way so the used column count
# returns to the workbook default (16384).
$ws.Range("Z1:AA1").EntireColumn.Delete()

# foreign_amount -> source_amount (same column, new meaning/label)
$ws.Range("B1").Value = "source_amount"

# Fill in the two newly inserted columns
$ws.Range("E1").Value = "target_amount"
$ws.Range("F1").Value = "target_fees"

# currency_conversions is now the sheet the user is looking at
$ws.Activate()
